# "Generate Report for Archive"
#
# The localization status changed from "Ready for handoff" to
# "In Translation" for the two handed-back files. That text lives in a
# single shared string that is reused by:
#   - Overview sheet : E2, F2 (zh-cn/de-de status for row 2) and E3, F3 (row 3)
#   - zh-cn sheet     : C2, C3 (Status column)
#   - de-de sheet     : C2, C3 (Status column)
# Updating every cell that displayed the old status text lets the engine
# dedupe them back down to one shared string entry, same as the diff shows.
#
# The shorter replacement text also narrowed the "Status" columns on all
# three sheets (Overview columns E/F, and column C on the two language
# sheets) - presumably from a column autosize pass when the report was
# regenerated.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
